# Scheduled runner update: refresh market-price derived figures in the
# per-job profit sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR). Values below
# mirror the latest Universalis price pull; BSM sheet is unchanged this run.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 71842.25
$ws.Range("I33").Value = 95445.19
$ws.Range("K33").Value = 95445.19
$ws.Range("M33").Value = -95216.19
$ws.Range("H51").Value = 8053.9653
$ws.Range("I51").Value = 8053.9653
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 8053.9653
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -7569.9653
$ws.Range("N51").ClearContents()
$ws.Range("H76").Value = 5200.5
$ws.Range("J76").Value = 5534.6665
$ws.Range("L76").Value = 5534.6665
$ws.Range("N76").Value = -6164.6665
$ws.Range("H79").Value = 5200.5
$ws.Range("J79").Value = 5534.6665
$ws.Range("L79").Value = 5534.6665
$ws.Range("N79").Value = -7718.6665
$ws.Range("H99").Value = 455
$ws.Range("I99").Value = 396.5
$ws.Range("J99").Value = 630.5
$ws.Range("K99").Value = 1189.5
$ws.Range("L99").Value = 1891.5
$ws.Range("M99").Value = 308.5
$ws.Range("N99").Value = -4887.5
$ws.Range("H101").Value = 19489.4
$ws.Range("J101").Value = 63331.332
$ws.Range("L101").Value = 189993.996
$ws.Range("N101").Value = -193237.996
$ws.Range("H116").Value = 3211
$ws.Range("J116").Value = 3275.625
$ws.Range("L116").Value = 3275.625
$ws.Range("N116").Value = -10159.625
$ws.Range("H137").Value = 1650.7273
$ws.Range("I137").Value = 1269
$ws.Range("J137").Value = 2108.8
$ws.Range("K137").Value = 3807
$ws.Range("L137").Value = 6326.400000000001
$ws.Range("M137").Value = -1257
$ws.Range("N137").Value = -11426.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6503.4
$ws.Range("I32").Value = 6730
$ws.Range("K32").Value = 6730
$ws.Range("M32").Value = -6443
$ws.Range("H45").Value = 2854.3125
$ws.Range("I45").Value = 839.5
$ws.Range("J45").Value = 4869.125
$ws.Range("K45").Value = 839.5
$ws.Range("L45").Value = 4869.125
$ws.Range("M45").Value = -462.5
$ws.Range("N45").Value = -5623.125
$ws.Range("H61").Value = 914.6667
$ws.Range("I61").Value = 914.6667
$ws.Range("K61").Value = 914.6667
$ws.Range("M61").Value = -702.6667
$ws.Range("H93").Value = 20389
$ws.Range("J93").Value = 20389
$ws.Range("L93").Value = 20389
$ws.Range("N93").Value = -25381
$ws.Range("H124").Value = 40812.145
$ws.Range("J124").Value = 40812.145
$ws.Range("L124").Value = 40812.145
$ws.Range("N124").Value = -50632.145
$ws.Range("H136").Value = 914.6667
$ws.Range("I136").Value = 914.6667
$ws.Range("K136").Value = 2744.0001
$ws.Range("M136").Value = -194.0001000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3997.5
$ws.Range("I86").Value = 3995
$ws.Range("K86").Value = 3995
$ws.Range("M86").Value = -2872
$ws.Range("H89").Value = 3997.5
$ws.Range("I89").Value = 3995
$ws.Range("K89").Value = 19975
$ws.Range("M89").Value = -14359
$ws.Range("H102").Value = 39999
$ws.Range("J102").Value = 39999
$ws.Range("L102").Value = 39999
$ws.Range("N102").Value = -44867
$ws.Range("H122").Value = 2699.6667
$ws.Range("I122").Value = 2699.6667
$ws.Range("K122").Value = 8099.000100000001
$ws.Range("M122").Value = -5649.000100000001
$ws.Range("H132").Value = 1597.25
$ws.Range("I132").Value = 1466.3334
$ws.Range("K132").Value = 4399.0002
$ws.Range("M132").Value = -1869.0002
$ws.Range("H134").Value = 2807.0732
$ws.Range("I134").Value = 2241.4412
$ws.Range("K134").Value = 6724.323600000001
$ws.Range("M134").Value = -4189.323600000001
$ws.Range("H141").Value = 184437.38
$ws.Range("I141").Value = 99000
$ws.Range("J141").Value = 190133.2
$ws.Range("K141").Value = 99000
$ws.Range("L141").Value = 190133.2
$ws.Range("M141").Value = -93820
$ws.Range("N141").Value = -200493.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2944.5
$ws.Range("I3").Value = 2944.5
$ws.Range("K3").Value = 8833.5
$ws.Range("M3").Value = -8721.5
$ws.Range("H4").Value = 168085.5
$ws.Range("I4").Value = 211.75
$ws.Range("K4").Value = 635.25
$ws.Range("M4").Value = -523.25
$ws.Range("H17").Value = 526.3077
$ws.Range("I17").Value = 154.25
$ws.Range("J17").Value = 691.6667
$ws.Range("K17").Value = 462.75
$ws.Range("L17").Value = 2075.0001
$ws.Range("M17").Value = -293.75
$ws.Range("N17").Value = -2413.0001
$ws.Range("H132").Value = 4796.1055
$ws.Range("I132").Value = 1486.9
$ws.Range("J132").Value = 8473
$ws.Range("K132").Value = 13382.1
$ws.Range("L132").Value = 76257
$ws.Range("M132").Value = -10852.1
$ws.Range("N132").Value = -81317

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2650
$ws.Range("J80").Value = 2687.625
$ws.Range("L80").Value = 2687.625
$ws.Range("N80").Value = -4683.625
$ws.Range("H83").Value = 2650
$ws.Range("J83").Value = 2687.625
$ws.Range("L83").Value = 13438.125
$ws.Range("N83").Value = -23422.125
$ws.Range("H132").Value = 3162.5
$ws.Range("I132").Value = 3126.25
$ws.Range("J132").Value = 3235
$ws.Range("K132").Value = 9378.75
$ws.Range("L132").Value = 9705
$ws.Range("M132").Value = -6848.75
$ws.Range("N132").Value = -14765

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2238.25
$ws.Range("I7").Value = 1967.3334
$ws.Range("J7").Value = 2509.1667
$ws.Range("K7").Value = 1967.3334
$ws.Range("L7").Value = 2509.1667
$ws.Range("M7").Value = -1855.3334
$ws.Range("N7").Value = -2733.1667
$ws.Range("H16").Value = 689.88
$ws.Range("I16").Value = 709.4583
$ws.Range("K16").Value = 709.4583
$ws.Range("M16").Value = -539.4583
$ws.Range("H46").Value = 2661.075
$ws.Range("J46").Value = 3575.1853
$ws.Range("L46").Value = 3575.1853
$ws.Range("N46").Value = -3951.1853
$ws.Range("H109").Value = 39799.5
$ws.Range("J109").Value = 39799.5
$ws.Range("L109").Value = 39799.5
$ws.Range("N109").Value = -42573.5
$ws.Range("H126").Value = 2238.25
$ws.Range("I126").Value = 1967.3334
$ws.Range("J126").Value = 2509.1667
$ws.Range("K126").Value = 5902.0002
$ws.Range("L126").Value = 7527.500100000001
$ws.Range("M126").Value = -3432.0002
$ws.Range("N126").Value = -12467.5001
$ws.Range("H132").Value = 12384.333
$ws.Range("I132").Value = 8999.4
$ws.Range("J132").Value = 16615.5
$ws.Range("K132").Value = 26998.2
$ws.Range("L132").Value = 49846.5
$ws.Range("M132").Value = -24468.2
$ws.Range("N132").Value = -54906.5
$ws.Range("H136").Value = 41670948
$ws.Range("I136").Value = 3373.2778
$ws.Range("K136").Value = 10119.8334
$ws.Range("M136").Value = -7569.8334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 463182.66
$ws.Range("J41").Value = 463182.66
$ws.Range("L41").Value = 463182.66
$ws.Range("N41").Value = -463962.66
$ws.Range("H62").Value = 4392.857
$ws.Range("J62").Value = 3500
$ws.Range("L62").Value = 3500
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 4392.857
$ws.Range("J65").Value = 3500
$ws.Range("L65").Value = 17500
$ws.Range("N65").Value = -23740
$ws.Range("H96").Value = 2156.8333
$ws.Range("I96").Value = 2889
$ws.Range("J96").Value = 692.5
$ws.Range("K96").Value = 2889
$ws.Range("L96").Value = 692.5
$ws.Range("M96").Value = -1516
$ws.Range("N96").Value = -3438.5
$ws.Range("H105").Value = 12999
$ws.Range("J105").Value = 12999
$ws.Range("L105").Value = 12999
$ws.Range("N105").Value = -19987
$ws.Range("H109").Value = 11568.25
$ws.Range("J109").Value = 11568.25
$ws.Range("L109").Value = 11568.25
$ws.Range("N109").Value = -14342.25

